# parser: add time front graph
# Adds a new summary table (rows 29-34) to the "295k" worksheet, mirroring
# the existing tables at rows 1-6 / 13-18 / 21-26, plus a new "golay filter"
# header cell in K29.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("295k")

# Header row (29)
$ws.Range("B29").Value = "all"
$ws.Range("C29").Value = "f1_good"
$ws.Range("D29").Value = "f1_bad"
$ws.Range("E29").Value = "f2_good"
$ws.Range("F29").Value = "f2_bad"
$ws.Range("G29").Value = "f3_good"
$ws.Range("H29").Value = "f3_bad"
$ws.Range("J29").Value = "chi2_per_dof_th"
$ws.Range("K29").Value = "golay filter"

# Totals row (30)
$ws.Range("A30").Value = "всего"
$ws.Range("B30").Value = 79
$ws.Range("D30").Value = 19
$ws.Range("F30").Value = 8
$ws.Range("H30").Value = 8
$ws.Range("J30").Value = 5

# шумы row (31) - label only
$ws.Range("A31").Value = "шумы"

# одиночные row (32)
$ws.Range("A32").Value = "одиночные"
$ws.Range("D32").Value = 2
$ws.Range("F32").Value = 2

# двойные row (33)
$ws.Range("A33").Value = "двойные"
$ws.Range("D33").Value = 11

# тройные row (34)
$ws.Range("A34").Value = "тройные"
$ws.Range("D34").Value = 6
$ws.Range("F34").Value = 6

# Scroll the view down to the new table and move the selection, matching
# the author's final cursor position after adding the table.
$ws.Activate()
$ws.Range("E36").Select()
